# Update countries & provincias Spain
#
# This script applies the 6-May-2020 13:06 data refresh:
#   - Updated totals for Estados Unidos and Libano
#   - Nepal's entry is (re)inserted just above Uganda with refreshed
#     counts; every country previously occupying rows 158-163 is pushed
#     down one row (each keeps its own totals, just shifted), and the
#     old stand-alone Nepal row (ex row 164) is absorbed by that shift.
#   - Liechtenstein/Libia row count refresh
#   - Montserrat is (re)placed just above Seychelles (the two rows swap
#     their country label + totals)
#   - Footer timestamp bumped from 12:33 to 13:06

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) -------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 6 de Mayo de 2020 a las 13:06"

# --- Estados Unidos (row 4): refreshed totals ------------------------------
$ws.Range("B4").Value2 = 1238052
$ws.Range("C4").Value2 = 419
$ws.Range("E4").Value2 = 965099

# --- Libano (row 106): refreshed totals ------------------------------------
$ws.Range("B106").Value2 = 750
$ws.Range("C106").Value2 = 9
$ws.Range("E106").Value2 = 519

# --- Nepal inserted right before Uganda; rows 158-164 shift down one -------
# Row 158: Nepal (brand new refreshed counts)
$ws.Range("A158").Value2 = "Nepal"
$ws.Range("B158").Value2 = 99
$ws.Range("C158").Value2 = 17
$ws.Range("D158").Value2 = 22
$ws.Range("E158").Value2 = 77
$ws.Range("F158").Value2 = 0
$ws.Range("G158").Value2 = 0
$ws.Range("H158").Value2 = 0

# Row 159: Uganda (formerly row 158's totals)
$ws.Range("A159").Value2 = "Uganda"
$ws.Range("B159").Value2 = 98
$ws.Range("C159").Value2 = 0
$ws.Range("D159").Value2 = 55
$ws.Range("E159").Value2 = 43
$ws.Range("F159").Value2 = 0
$ws.Range("G159").Value2 = 0
$ws.Range("H159").Value2 = 0

# Row 160: Benin (formerly row 159's totals)
$ws.Range("A160").Value2 = "Benin"
$ws.Range("B160").Value2 = 96
$ws.Range("C160").Value2 = 0
$ws.Range("D160").Value2 = 50
$ws.Range("E160").Value2 = 44
$ws.Range("F160").Value2 = 0
$ws.Range("G160").Value2 = 0
$ws.Range("H160").Value2 = 2

# Row 161: Monaco (formerly row 160's totals)
$ws.Range("A161").Value2 = "Monaco"
$ws.Range("B161").Value2 = 95
$ws.Range("C161").Value2 = 0
$ws.Range("D161").Value2 = 81
$ws.Range("E161").Value2 = 10
$ws.Range("F161").Value2 = 1
$ws.Range("G161").Value2 = 0
$ws.Range("H161").Value2 = 4

# Row 162: Guyana (formerly row 161's totals)
$ws.Range("A162").Value2 = "Guyana"
$ws.Range("B162").Value2 = 93
$ws.Range("C162").Value2 = 1
$ws.Range("D162").Value2 = 27
$ws.Range("E162").Value2 = 56
$ws.Range("F162").Value2 = 3
$ws.Range("G162").Value2 = 1
$ws.Range("H162").Value2 = 10

# Row 163: Bahamas (formerly row 162's totals)
$ws.Range("A163").Value2 = "Bahamas"
$ws.Range("B163").Value2 = 89
$ws.Range("C163").Value2 = 0
$ws.Range("D163").Value2 = 26
$ws.Range("E163").Value2 = 52
$ws.Range("F163").Value2 = 1
$ws.Range("G163").Value2 = 0
$ws.Range("H163").Value2 = 11

# Row 164: Republica de Africa Central (formerly row 163's totals); this
# also absorbs the slot that used to hold the stand-alone Nepal row.
$ws.Range("A164").Value2 = "Republica de Africa Central"
$ws.Range("B164").Value2 = 85
$ws.Range("C164").Value2 = 0
$ws.Range("D164").Value2 = 10
$ws.Range("E164").Value2 = 75
$ws.Range("F164").Value2 = 0
$ws.Range("G164").Value2 = 0
$ws.Range("H164").Value2 = 0

# --- Libia (row 170): refreshed totals -------------------------------------
$ws.Range("D170").Value2 = 24
$ws.Range("E170").Value2 = 36

# --- Montserrat swaps above Seychelles (rows 205/206) -----------------------
# Row 205: Montserrat (formerly Seychelles' row, now Montserrat's totals)
$ws.Range("A205").Value2 = "Montserrat"
$ws.Range("B205").Value2 = 11
$ws.Range("C205").Value2 = 0
$ws.Range("D205").Value2 = 7
$ws.Range("E205").Value2 = 3
$ws.Range("F205").Value2 = 1
$ws.Range("G205").Value2 = 0
$ws.Range("H205").Value2 = 1

# Row 206: Seychelles (formerly Montserrat's row, now Seychelles' totals)
$ws.Range("A206").Value2 = "Seychelles"
$ws.Range("B206").Value2 = 11
$ws.Range("C206").Value2 = 0
$ws.Range("D206").Value2 = 8
$ws.Range("E206").Value2 = 3
$ws.Range("F206").Value2 = 0
$ws.Range("G206").Value2 = 0
$ws.Range("H206").Value2 = 0
